# Swap the two embedded themes:
#   ppt/theme/theme1.xml  ("Integral", used by the slide master)   -> "Office Theme" palette
#   ppt/theme/theme2.xml  ("Office Theme", used by the notes master) -> "Integral" palette
#
# The PowerPoint object model doesn't give direct access to the raw OOXML
# theme parts, so we reproduce the swap through the Design/Theme color
# scheme API: every ThemeColorScheme slot on the slide master's Theme is
# set to the RGB value the target "Office Theme" palette uses. (dk1/lt1 -
# black/white - are identical in both themes already, so only the other
# ten slots actually need new values.)

$p  = $ppt.ActivePresentation
$sm = $p.SlideMaster
$t  = $sm.Theme
$cs = $t.ThemeColorScheme

function Set-ThemeRgb($scheme, $index, $hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    $bgr = ($b * 65536) + ($g * 256) + $r
    $scheme.Item($index).RGB = $bgr
}

# Office Theme color scheme (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink)
Set-ThemeRgb $cs 1  "000000"
Set-ThemeRgb $cs 2  "FFFFFF"
Set-ThemeRgb $cs 3  "44546A"
Set-ThemeRgb $cs 4  "E7E6E6"
Set-ThemeRgb $cs 5  "5B9BD5"
Set-ThemeRgb $cs 6  "ED7D31"
Set-ThemeRgb $cs 7  "A5A5A5"
Set-ThemeRgb $cs 8  "FFC000"
Set-ThemeRgb $cs 9  "4472C4"
Set-ThemeRgb $cs 10 "70AD47"
Set-ThemeRgb $cs 11 "0563C1"
Set-ThemeRgb $cs 12 "954F72"

# Best-effort: rename the theme/color-scheme to match (the host may not
# expose a writable Name on these objects; ignore failures silently).
try { $t.Name = "Office Theme" } catch {}
try { $cs.Name = "Office" } catch {}
